$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column widths (only columns whose width changed per diff) ---
$ws.Columns.Item(2).ColumnWidth = 7.166666666666667  # B -> width 8
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667  # C -> width 8
$ws.Columns.Item(5).ColumnWidth = 7.166666666666667  # E -> width 8
$ws.Columns.Item(6).ColumnWidth = 7.166666666666667  # F -> width 8
$ws.Columns.Item(7).ColumnWidth = 7.166666666666667  # G -> width 8
$ws.Columns.Item(9).ColumnWidth = 7.166666666666667  # I -> width 8
$ws.Columns.Item(10).ColumnWidth = 7.166666666666667  # J -> width 8
$ws.Columns.Item(11).ColumnWidth = 7.166666666666667  # K -> width 8
$ws.Columns.Item(12).ColumnWidth = 7.166666666666667  # L -> width 8
$ws.Columns.Item(13).ColumnWidth = 7.166666666666667  # M -> width 8
$ws.Columns.Item(15).ColumnWidth = 7.166666666666667  # O -> width 8
$ws.Columns.Item(16).ColumnWidth = 7.166666666666667  # P -> width 8
$ws.Columns.Item(17).ColumnWidth = 7.166666666666667  # Q -> width 8
$ws.Columns.Item(20).ColumnWidth = 8.166666666666666  # T -> width 9
$ws.Columns.Item(21).ColumnWidth = 7.166666666666667  # U -> width 8
$ws.Columns.Item(22).ColumnWidth = 7.166666666666667  # V -> width 8
$ws.Columns.Item(23).ColumnWidth = 7.166666666666667  # W -> width 8
$ws.Columns.Item(24).ColumnWidth = 7.166666666666667  # X -> width 8
$ws.Columns.Item(26).ColumnWidth = 7.166666666666667  # Z -> width 8
$ws.Columns.Item(27).ColumnWidth = 7.166666666666667  # AA -> width 8
$ws.Columns.Item(28).ColumnWidth = 7.166666666666667  # AB -> width 8
$ws.Columns.Item(29).ColumnWidth = 7.166666666666667  # AC -> width 8
$ws.Columns.Item(30).ColumnWidth = 7.166666666666667  # AD -> width 8
$ws.Columns.Item(34).ColumnWidth = 7.166666666666667  # AH -> width 8

# --- Update cell values for rows 2-5 with the new dataset ---
# Row 2
$ws.Range("A2").Value = 45121.50694444445
$ws.Range("B2").Value = 24.021
$ws.Range("C2").Value = 16.687
$ws.Range("D2").Value = 4.255
$ws.Range("E2").Value = 50.696
$ws.Range("F2").Value = 41.946
$ws.Range("G2").Value = 18.904
$ws.Range("H2").Value = 63.181
$ws.Range("I2").Value = 29.086
$ws.Range("J2").Value = 12.432
$ws.Range("K2").Value = 19.179
$ws.Range("L2").Value = 19.925
$ws.Range("M2").Value = 20.923
$ws.Range("N2").Value = 6.036
$ws.Range("O2").Value = 18.798
$ws.Range("P2").Value = 26.508
$ws.Range("Q2").Value = 15.598
$ws.Range("R2").Value = 3.832
$ws.Range("S2").Value = 2.606
$ws.Range("T2").Value = 278.79
$ws.Range("U2").Value = 52.368
$ws.Range("V2").Value = 17.351
$ws.Range("W2").Value = 34.892
$ws.Range("X2").Value = 18.126
$ws.Range("Y2").Value = 2.393
$ws.Range("Z2").Value = 31.289
$ws.Range("AA2").Value = 15.326
$ws.Range("AB2").Value = 13.706
$ws.Range("AC2").Value = 16.029
$ws.Range("AD2").Value = 20.711
$ws.Range("AE2").Value = 3.641
$ws.Range("AF2").Value = 55.941
$ws.Range("AG2").Value = 9.720000000000001
$ws.Range("AH2").Value = 21.693

# Row 3
$ws.Range("A3").Value = 45121.51388888889
$ws.Range("B3").Value = 11.53
$ws.Range("C3").Value = 7.99
$ws.Range("D3").Value = 1.651
$ws.Range("E3").Value = 24.472
$ws.Range("F3").Value = 20.288
$ws.Range("G3").Value = 9.074
$ws.Range("H3").Value = 38.153
$ws.Range("I3").Value = 13.961
$ws.Range("J3").Value = 5.973
$ws.Range("K3").Value = 9.134
$ws.Range("L3").Value = 9.702
$ws.Range("M3").Value = 10.195
$ws.Range("N3").Value = 2.901
$ws.Range("O3").Value = 9.023
$ws.Range("P3").Value = 12.709
$ws.Range("Q3").Value = 7.741
$ws.Range("R3").Value = 1.596
$ws.Range("S3").Value = 0.922
$ws.Range("T3").Value = 130.034
$ws.Range("U3").Value = 25.359
$ws.Range("V3").Value = 8.329000000000001
$ws.Range("W3").Value = 16.769
$ws.Range("X3").Value = 8.959
$ws.Range("Y3").Value = 1.119
$ws.Range("Z3").Value = 17.834
$ws.Range("AA3").Value = 7.357
$ws.Range("AB3").Value = 6.705
$ws.Range("AC3").Value = 7.833
$ws.Range("AD3").Value = 10.152
$ws.Range("AE3").Value = 1.294
$ws.Range("AF3").Value = 34.546
$ws.Range("AG3").Value = 4.611
$ws.Range("AH3").Value = 10.413

# Row 4
$ws.Range("A4").Value = 45121.52083333334
$ws.Range("B4").Value = 1.922
$ws.Range("C4").Value = 0.998
$ws.Range("D4").Value = 0.838
$ws.Range("E4").Value = 3.782
$ws.Range("F4").Value = 3.179
$ws.Range("G4").Value = 1.515
$ws.Range("H4").Value = 11.632
$ws.Range("I4").Value = 2.327
$ws.Range("J4").Value = 0.9399999999999999
$ws.Range("K4").Value = 1.38
$ws.Range("L4").Value = 1.473
$ws.Range("M4").Value = 1.494
$ws.Range("N4").Value = 0.505
$ws.Range("O4").Value = 1.504
$ws.Range("P4").Value = 2.131
$ws.Range("Q4").Value = 1.494
$ws.Range("R4").Value = 0.968
$ws.Range("S4").Value = 0.383
$ws.Range("T4").Value = 15.66
$ws.Range("U4").Value = 4.6
$ws.Range("V4").Value = 1.388
$ws.Range("W4").Value = 2.937
$ws.Range("X4").Value = 1.642
$ws.Range("Y4").Value = 0.139
$ws.Range("Z4").Value = 4.884
$ws.Range("AA4").Value = 1.226
$ws.Range("AB4").Value = 1.271
$ws.Range("AC4").Value = 1.448
$ws.Range("AD4").Value = 1.54
$ws.Range("AE4").Value = 0.784
$ws.Range("AF4").Value = 10.757
$ws.Range("AG4").Value = 0.676
$ws.Range("AH4").Value = 1.743

# Row 5
$ws.Range("A5").Value = 45121.52777777778
$ws.Range("B5").Value = 1.44
$ws.Range("C5").Value = 0.75
$ws.Range("D5").Value = 0.6
$ws.Range("E5").Value = 2.84
$ws.Range("F5").Value = 2.39
$ws.Range("G5").Value = 1.13
$ws.Range("H5").Value = 6.62
$ws.Range("I5").Value = 1.75
$ws.Range("J5").Value = 0.64
$ws.Range("K5").Value = 1.02
$ws.Range("L5").Value = 1.12
$ws.Range("M5").Value = 1.13
$ws.Range("N5").Value = 0.37
$ws.Range("O5").Value = 1.13
$ws.Range("P5").Value = 1.51
$ws.Range("Q5").Value = 1.15
$ws.Range("R5").Value = 0.71
$ws.Range("S5").Value = 0.27
$ws.Range("T5").Value = 9.84
$ws.Range("U5").Value = 3.25
$ws.Range("V5").Value = 1.04
$ws.Range("W5").Value = 1.96
$ws.Range("X5").Value = 1.2
$ws.Range("Y5").Value = 0.09
$ws.Range("Z5").Value = 2.69
$ws.Range("AA5").Value = 0.92
$ws.Range("AB5").Value = 0.96
$ws.Range("AC5").Value = 1.09
$ws.Range("AD5").Value = 1.18
$ws.Range("AE5").Value = 0.5600000000000001
$ws.Range("AF5").Value = 5.91
$ws.Range("AG5").Value = 0.5
$ws.Range("AH5").Value = 1.3

# --- Remove old row 6 (dataset reduced from 5 to 4 data rows) ---
$ws.Rows.Item(6).Delete()
